$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = "Overall Average Throughput:"
$ws.Range("H3").Formula = "=AVERAGE(E2:E261)"

$ws.Range("G5").Value = "Pipes"
$ws.Range("H5").Value = "Average Throughput"
$ws.Range("I5").Value = "Highest Throughput"

$ws.Range("G6").Value = 1
$ws.Range("H6").Formula = "=AVERAGE(E2:E41)"
$ws.Range("I6").Formula = "=MAX(E2:E41)"

$ws.Range("G7").Value = 2
$ws.Range("H7").Formula = "=AVERAGE(E42:E81)"
$ws.Range("I7").Formula = "=MAX(E42:E81)"

$ws.Range("G8").Value = 3
$ws.Range("H8").Formula = "=AVERAGE(E82:E121)"
$ws.Range("I8").Formula = "=MAX(E82:E121)"

$ws.Range("G9").Value = 4
$ws.Range("H9").Formula = "=AVERAGE(E122:E161)"
$ws.Range("I9").Formula = "=MAX(E122:E161)"

$ws.Range("G10").Value = 5
$ws.Range("H10").Formula = "=AVERAGE(E162:E201)"
$ws.Range("I10").Formula = "=MAX(E162:E201)"

$ws.Range("G11").Value = 10
$ws.Range("H11").Formula = "=AVERAGE(E202:E231)"
$ws.Range("I11").Formula = "=MAX(E202:E231)"

$ws.Range("G12").Value = 16
$ws.Range("H12").Formula = "=AVERAGE(E232:E261)"
$ws.Range("I12").Formula = "=MAX(E232:E261)"

$ws.Range("A1").Copy()
$ws.Range("G5:I5").PasteSpecial(-4122)

$ws.Range("I13").Select()
$ws.Application.ActiveWindow.ScrollColumn = 6
